$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 701, shifting existing rows 701:750 down to 702:751
$ws.Rows("701:701").Insert()

# Populate the new row 701 with its values
$ws.Range("A701").Value = 6
$ws.Range("B701").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C701").Value = "Metropolitana"
$ws.Range("D701").Value = 44753
$ws.Range("E701").Value = 13
$ws.Range("F701").Value = 100112003
$ws.Range("G701").Value = "Ajo"
$ws.Range("H701").Value = "Chino"
$ws.Range("I701").Value = "Primera"
$ws.Range("J701").Value = 1100
$ws.Range("K701").Value = 21000
$ws.Range("L701").Value = 22000
$ws.Range("M701").Value = 21545
$ws.Range("N701").Value = "`$/caja 10 kilos"
$ws.Range("O701").Value = "China"
$ws.Range("P701").Value = 2154
$ws.Range("Q701").Value = 10
$ws.Range("R701").Value = "Hortaliza"
